$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 0.09056462818037136
$ws.Range("H2").Value = 0.2126375042666278
$ws.Range("K2").Value = 3.892599562066422
$ws.Range("L2").Value = "[-0.9825385831595881, 8.767737707292431]"
$ws.Range("M2").Value = 0.116920322616666
$ws.Range("N2").Value = 0.116920322616666
$ws.Range("O2").Value = -1.01889491467777
$ws.Range("P2").Value = "[-3.610158524845928, 1.572368695490387]"
$ws.Range("Q2").Value = 0.4389357108471779
$ws.Range("R2").Value = 0.4389357108471779
$ws.Range("S2").Value = 14.82349409064097
$ws.Range("T2").Value = "[12.326124570009817, 17.32086361127213]"
$ws.Range("W2").Value = 4.214594594594693
$ws.Range("X2").Value = -6.504004004004162
$ws.Range("Y2").Value = 14.93319319319355

# Row 3 updates
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 22.86000000000013
$ws.Range("G3").Value = 0.007304566710637284
$ws.Range("H3").Value = 0.05751643789443971
$ws.Range("I3").Value = ""
$ws.Range("K3").Value = 4.429394970564768
$ws.Range("L3").Value = "[0.9705874601093161, 7.888202481020221]"
$ws.Range("M3").Value = 0.0122969344897601
$ws.Range("N3").Value = 0.02459386897952021
$ws.Range("O3").Value = 1.83023716155081
$ws.Range("P3").Value = "[0.823921196436963, 2.8365531266646578]"
$ws.Range("Q3").Value = 0.0004122375303228232
$ws.Range("R3").Value = 0.0008244750606456464
$ws.Range("S3").Value = 12.90817195994105
$ws.Range("T3").Value = "[10.953129788742793, 14.863214131139301]"
$ws.Range("W3").Value = 16.20108108108118
$ws.Range("X3").Value = 12.5398198198199
$ws.Range("Y3").Value = 19.86234234234246
